$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.098742
$ws.Range("H2").Value = 0.296226
$ws.Range("M2").Value = 8.333446333333333
$ws.Range("N2").Value = 25.000339
$ws.Range("O2").Value = 0.3294294409523786
$ws.Range("P2").Value = 0.3294294409523787
$ws.Range("Q2").Value = 0.8228611578459999
$ws.Range("R2").Value = 7.405750420614
$ws.Range("S2").Value = 0.3294294409523786
$ws.Range("T2").Value = 0.3294294409523787

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.098742
$ws.Range("H3").Value = 0.296226
$ws.Range("O3").Value = 0.357886883212021
$ws.Range("P3").Value = 0.357886883212021
$ws.Range("Q3").Value = 0.8939432196659999
$ws.Range("R3").Value = 8.045488976993999
$ws.Range("S3").Value = 0.357886883212021
$ws.Range("T3").Value = 0.357886883212021

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.098742
$ws.Range("H4").Value = 0.296226
$ws.Range("M4").Value = 4.309709000000001
$ws.Range("N4").Value = 12.929127
$ws.Range("O4").Value = 0.170367093006711
$ws.Range("P4").Value = 0.170367093006711
$ws.Range("Q4").Value = 0.4255492860780001
$ws.Range("R4").Value = 3.829943574702
$ws.Range("S4").Value = 0.170367093006711
$ws.Range("T4").Value = 0.170367093006711

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.098742
$ws.Range("H5").Value = 0.296226
$ws.Range("M5").Value = 3.600126333333333
$ws.Range("N5").Value = 10.800379
$ws.Range("O5").Value = 0.1423165828288893
$ws.Range("P5").Value = 0.1423165828288893
$ws.Range("Q5").Value = 0.355483674406
$ws.Range("R5").Value = 3.199353069654
$ws.Range("S5").Value = 0.1423165828288893
$ws.Range("T5").Value = 0.1423165828288893

$wb.Save()
